# "add rural communities interactions back in, fix de_dg files (not run yet)"
#
# The "Legislature" stakeholder row (row 12) is no longer needed on this
# sheet, so the whole row is deleted. Excel's Delete() on a full row shifts
# every row below it up by one (values, styles, row heights all move with
# it), which is exactly what the target workbook shows - every row from 13
# downward moved up to become row 12 downward, and the final row (29) is
# gone because there's nothing left to shift into it.
#
# Deleting the row also drops "Legislature" from the shared-string table,
# which is why every shared-string index above it shifts down by one in the
# saved file - that happens automatically when the workbook is saved, we
# just need to remove the row/cells that referenced it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(12).Delete()

# Leave the selection where the user's edit left it - on the row that moved
# up into the deleted row's place (now occupied by "Friant Water
# Authority"), selected as a full row just like after a row-delete in the UI.
$ws.Range("A12:XFD12").Select()
